# Auto-generated edit script: updates cryptos price/volume table
# Mirrors a scheduled GitHub Actions data refresh of the cryptos list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D hold plain numeric-looking text (e.g. grouped
# thousands like '30.578.19' or zero-padded decimals like '0.07700').
# Force Text format before writing so Excel keeps the literal digits
# instead of auto-converting them into a Number/Date value.
$priceCells = @('D2','D3','D5','D7','D8','D9','D10','D11','D12','D14','D15','D16','D17','D18','D20','D21','D23','D24','D25','D26','D27','D28','D29','D32','D33','D34','D38','D39','D40','D41','D42','D43','D45','D46','D47','D48','D49','D50','D51')
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.578.19'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '1.874.12'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '248.33'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4729'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '0.2911'
$ws.Range('D9').Value = '0.06479'
$ws.Range('D10').Value = '22.11'
$ws.Range('E10').Value = '  +5.19%  '
$ws.Range('D11').Value = '0.07700'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').Value = '0.7415'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').Value = '1.873.89'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '5.171'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = '273.82'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '30.625.70'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '13.33'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '0.000007509'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = '2.115.36'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '5.264'
$ws.Range('E23').Value = '  +0.64%  '
$ws.Range('D24').Value = '6.186'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '9.209'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').Value = '164.41'
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').Value = '18.74'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('D28').Value = '1.910'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '0.1003'
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').Value = '4.271'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').Value = '4.097'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '0.04795'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '0.01853'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '2.754'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').Value = '6.234'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('D41').Value = '73.33'
$ws.Range('E41').Value = '  +5.17%  '
$ws.Range('D42').Value = '1.970'
$ws.Range('E42').Value = '  +3.12%  '
$ws.Range('D43').Value = '0.4175'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '0.8348'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').Value = '102.03'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').Value = '9.322'
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '7.009'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '35.44'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').Value = '915.72'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').Value = '0.05655'
$ws.Range('E51').Value = '  +1.54%  '
